# "Generate Report for Archive"
#
# 1) Status text: "Ready for handoff" -> "In Translation" on every sheet
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2) Narrow the "Latest HO Xliff Generate Date" status columns:
#    Overview!E:F and the "Status" column (C) on the zh-cn / de-de sheets,
#    from ~17.22 chars down to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# --- 1) Replace status text everywhere it appears ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- 2) Resize the status-date columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C
